$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the tiny floating point precision of A17's timestamp
$ws.Cells.Item(17, 1).Value = 45876.62519435185

# Add new row 18
$ws.Cells.Item(18, 1).Value = 45876.66685023394
$ws.Cells.Item(18, 2).Value = 2025
$ws.Cells.Item(18, 3).Value = 28
$ws.Cells.Item(18, 4).Value = 16.78
$ws.Cells.Item(18, 5).Value = 84.06
$ws.Cells.Item(18, 6).Value = 166.99
$ws.Cells.Item(18, 7).Value = 14.68
$ws.Cells.Item(18, 8).Value = "ESE"
$ws.Cells.Item(18, 9).Value = 0
$ws.Cells.Item(18, 10).Value = "16:00:15"

# Match the date/time number format used by column A elsewhere
$ws.Cells.Item(18, 1).NumberFormat = $ws.Cells.Item(17, 1).NumberFormat
